$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '42.821.69'
$ws.Range("E2").Value = '  -2.30%  '

$ws.Range("D3").Value = '2.239.23'
$ws.Range("E3").Value = '  -2.24%  '

$cD = $ws.Range("D4")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '1.00'
$cD.Style = $origStyle
$ws.Range("E4").Value = '  +0.14%  '

$cD = $ws.Range("D5")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '112.66'
$cD.Style = $origStyle
$ws.Range("E5").Value = '  -5.94%  '

$cD = $ws.Range("D6")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '298.35'
$cD.Style = $origStyle
$ws.Range("E6").Value = '  +11.42%  '

$cD = $ws.Range("D7")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '0.628'
$cD.Style = $origStyle
$ws.Range("E7").Value = '  -2.05%  '

$ws.Range("E8").Value = '  +0.02%  '

$cD = $ws.Range("D9")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '0.615'
$cD.Style = $origStyle
$ws.Range("E9").Value = '  -0.47%  '

$cD = $ws.Range("D10")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '45.41'
$cD.Style = $origStyle
$ws.Range("E10").Value = '  -6.28%  '

$cD = $ws.Range("D11")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '0.0928'
$cD.Style = $origStyle
$ws.Range("E11").Value = '  -1.51%  '

$cD = $ws.Range("D12")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '55.80'
$cD.Style = $origStyle
$ws.Range("E12").Value = '  +2.15%  '

$cD = $ws.Range("D13")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '9.05'
$cD.Style = $origStyle
$ws.Range("E13").Value = '  -2.48%  '

$ws.Range("E14").Value = '  -2.97%  '

$cD = $ws.Range("D15")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '15.27'
$cD.Style = $origStyle
$ws.Range("E15").Value = '  -2.10%  '

$cD = $ws.Range("D16")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '0.893'
$cD.Style = $origStyle
$ws.Range("E16").Value = '  +0.04%  '

$ws.Range("D17").Value = '2.577.42'
$ws.Range("E17").Value = '  -2.13%  '

$ws.Range("D18").Value = '2.241.52'
$ws.Range("E18").Value = '  -1.91%  '

$ws.Range("D19").Value = '42.644.47'
$ws.Range("E19").Value = '  -2.45%  '

$cD = $ws.Range("D20")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '7.58'
$cD.Style = $origStyle
$ws.Range("E20").Value = '  +7.59%  '

$cD = $ws.Range("D21")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '0.0000108'
$cD.Style = $origStyle
$ws.Range("E21").Value = '  -2.10%  '

$cD = $ws.Range("D22")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '73.26'
$cD.Style = $origStyle
$ws.Range("E22").Value = '  +1.09%  '

$ws.Range("E23").Value = '  +23.51%  '

$cD = $ws.Range("D24")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '2.33'
$cD.Style = $origStyle
$ws.Range("E24").Value = '  -5.78%  '

$cD = $ws.Range("D25")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '231.21'
$cD.Style = $origStyle
$ws.Range("E25").Value = '  -2.04%  '

$cD = $ws.Range("D26")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '9.38'
$cD.Style = $origStyle
$ws.Range("E26").Value = '  -3.13%  '

$cD = $ws.Range("D27")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '12.02'
$cD.Style = $origStyle
$ws.Range("E27").Value = '  +1.38%  '

$ws.Range("E28").Value = '  -1.36%  '

$cD = $ws.Range("D29")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '39.73'
$cD.Style = $origStyle
$ws.Range("E29").Value = '  -8.28%  '

$cD = $ws.Range("D30")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '2.24'
$cD.Style = $origStyle
$ws.Range("E30").Value = '  -0.29%  '

$ws.Range("E31").Value = '  -3.94%  '

$cD = $ws.Range("D32")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '173.97'
$cD.Style = $origStyle
$ws.Range("E32").Value = '  +0.24%  '

$cD = $ws.Range("D33")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '21.25'
$cD.Style = $origStyle
$ws.Range("E33").Value = '  -2.42%  '

$cD = $ws.Range("D34")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '0.0894'
$cD.Style = $origStyle
$ws.Range("E34").Value = '  -2.20%  '

$cD = $ws.Range("D35")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '5.71'
$cD.Style = $origStyle
$ws.Range("E35").Value = '  -1.04%  '

$cD = $ws.Range("D36")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '4.97'
$cD.Style = $origStyle
$ws.Range("E36").Value = '  +4.86%  '

$cD = $ws.Range("D37")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '4.35'
$cD.Style = $origStyle
$ws.Range("E37").Value = '  +10.05%  '

$cD = $ws.Range("D38")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '0.128'
$cD.Style = $origStyle
$ws.Range("E38").Value = '  -1.60%  '

$cD = $ws.Range("D39")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '0.0371'
$cD.Style = $origStyle
$ws.Range("E39").Value = '  -3.12%  '

$ws.Range("E40").Value = '  -1.81%  '

$cD = $ws.Range("D41")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '2.56'
$cD.Style = $origStyle
$ws.Range("E41").Value = '  -0.40%  '

$cD = $ws.Range("D42")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '0.238'
$cD.Style = $origStyle
$ws.Range("E42").Value = '  -0.38%  '

$cD = $ws.Range("D43")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '72.13'
$cD.Style = $origStyle
$ws.Range("E43").Value = '  -3.93%  '

$cD = $ws.Range("D44")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '13.20'
$cD.Style = $origStyle
$ws.Range("E44").Value = '  -8.13%  '

$cD = $ws.Range("D45")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '1.00'
$cD.Style = $origStyle
$ws.Range("E45").Value = '  +0.34%  '

$ws.Range("E46").Value = '  -3.25%  '

$cD = $ws.Range("D47")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '5.59'
$cD.Style = $origStyle
$ws.Range("E47").Value = '  -6.55%  '

$cD = $ws.Range("D48")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '1.32'
$cD.Style = $origStyle
$ws.Range("E48").Value = '  +3.11%  '

$cD = $ws.Range("D49")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '105.63'
$cD.Style = $origStyle
$ws.Range("E49").Value = '  +3.53%  '

$cD = $ws.Range("D50")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '8.64'
$cD.Style = $origStyle
$ws.Range("E50").Value = '  +0.22%  '

$cD = $ws.Range("D51")
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = '0.0986'
$cD.Style = $origStyle
$ws.Range("E51").Value = '  -2.24%  '
